{"js": "// Objective paragraph update:\n//  1. \"about 3 years\" -> \"3+ years\"\n//  2. \"With the knowledge of \" -> \"With the knowledge of competitive programming and \"\nconst yearsPhrase = context.document.body.search(\"about 3 years\", { matchCase: true });\nyearsPhrase.load(\"text\");\n\nconst knowledgePhrase = context.document.body.search(\"With the knowledge of \", { matchCase: true });\nknowledgePhrase.load(\"text\");\n\nawait context.sync();\n\nif (yearsPhrase.items.length > 0) {\n  yearsPhrase.items[0].insertText(\"3+ years\", Word.InsertLocation.replace);\n}\n\nif (knowledgePhrase.items.length > 0) {\n  knowledgePhrase.items[0].insertText(\n    \"With the knowledge of competitive programming and \",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Objective paragraph update:\n#  1. \"about 3 years\" -> \"3+ years\"\n#  2. \"With the knowledge of \" -> \"With the knowledge of competitive programming and \"\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"about 3 years\", $false, $false, $false, $false, $false, $true, 1, $false, \"3+ years\", 2)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"With the knowledge of \", $false, $false, $false, $false, $false, $true, 1, $false, \"With the knowledge of competitive programming and \", 2)\n"}
